$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly added existing VRES and BESS capacities
$ws.Range("E7").Value = 12
$ws.Range("S7").Value = 8
$ws.Range("S8").Value = 8
$ws.Range("E9").Value = 39
$ws.Range("S9").Value = 8
$ws.Range("S10").Value = 8
$ws.Range("S11").Value = 8

# Update the frozen-pane selection to the new working cell
$ws.Range("J20").Select()
